$d = $word.ActiveDocument

# --- 1. Mark the run holding the screenshot drawing as NoProof --------------
# (adds <w:rPr><w:noProof/></w:rPr> to the run that contains the <w:drawing>)
$imgPara = $d.Paragraphs(2)
$imgPara.Range.NoProofing = $true

# --- 2. Grab a copy (incl. paragraph formatting/numbering) of the first ----
#        bullet paragraph so the new bullet reuses the same list (numId).
$firstPara = $d.Paragraphs(1)
$bulletSrc = $d.Range($firstPara.Range.Start, $firstPara.Range.End)
$bulletFormattedText = $bulletSrc.FormattedText

# --- 3. Append two blank paragraphs after the picture paragraph, then the --
#        new bulleted finding, reusing the ListParagraph style + numId 1.
$tail = $d.Range($d.Content.End, $d.Content.End)
$tail.InsertParagraphAfter()
$tail = $d.Range($d.Content.End, $d.Content.End)
$tail.InsertParagraphAfter()
$tail = $d.Range($d.Content.End, $d.Content.End)
$tail.InsertParagraphAfter()

$pasteRange = $d.Range($d.Content.End, $d.Content.End)
$pasteRange.FormattedText = $bulletFormattedText

# Pasting the formatted paragraph leaves one spare trailing empty paragraph
# behind it (Word always keeps a final paragraph mark) - drop it.
$d.Paragraphs.Last.Range.Delete()

# Swap in the real text for the newly added bullet paragraph.
$newBullet = $d.Paragraphs.Last
$newBullet.Range.Text = "Some patients which have multiple lesitions can have both cancerous and non concerous"
